$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted as row 46, pushing the existing
# rows 46-125 down to 47-126. Nothing else on the sheet changes - it is
# purely a one-row insert followed by filling in the new record's data.
$ws.Rows.Item(46).Insert()

# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T carry the same constant metadata as
# every other "Mango / Vega Monumental Concepción" row in the sheet.
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"

# Column D holds dates and uses the sheet's date number format.
$ws.Cells.Item(46, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 4).Value = 44804

$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100108
$ws.Cells.Item(46, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(46, 9).Value = 100108002
$ws.Cells.Item(46, 10).Value = "Mango"
$ws.Cells.Item(46, 11).Value = "Sin especificar"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 200
$ws.Cells.Item(46, 14).Value = 9000
$ws.Cells.Item(46, 15).Value = 9500
$ws.Cells.Item(46, 16).Value = 9250
$ws.Cells.Item(46, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(46, 18).Value = "Brasil"
$ws.Cells.Item(46, 19).Value = 2312
$ws.Cells.Item(46, 20).Value = 4
